# Daily update of covid19 tracker data files
# - Update NZL fiscal/monetary policy text (H43:K43)
# - Update TUR travel-ban and health text (D54, G54)
# - Bump "Updated on" date for all existing rows (43923 -> 43924)
# - Append 17 new blank rows (58-74) carrying only the updated date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

# --- NZL: Fiscal measures - overall (H43) ---
$h43 = @'
The government announced a massive package of support measures on 17 March amounting to NZD 12.1 billion (4% of GDP), with one half to be implemented in the coming year. Most (NZD 8.7 billion) of the fiscal measures are to support businesses and jobs while NZD 2.8 billion were allocated to increase social benefit payments and NZD 500 million (0.2% of GDP) were allocated to the health sector. On 27 March, the government revamped the fiscal package by expanding the coverage of the wage subsidy scheme. The size of fiscal package can be up to NZD 19 billlion (6.1% of GDP) depending on business uptake of wage subsidy. This package comes on top of the NZD 12 billion increase in infrastructure spending announced in January.     
'@

$ws.Range("H43").Value = $h43

# --- NZL: Fiscal measures - people specific (I43) ---
$i43 = @'
On 27 March, the government expanded coverage of the wage subsidy scheme. It is now expected to cost NZD 8-12 billion depending on business uptake.
'The NZD 2.8 billion increase in benefit payments will come in the form of a NZD 25 per week increase in core benefit payments and a doubling in the winter energy payment to $1400 for couples and $900 for single people. Also, the threshold for In Work Tax Credit has been removed.
A NZD 27 million package is being provided to social sector services and community groups to bolster existing essential social sector services,provide support for disabled people in lock-down, and fund innovative community-led solutions to support local resilience.  
'@

$ws.Range("I43").Value = $i43

# --- NZL: Fiscal measures - company specific (J43) ---
$j43 = @'
Firms are to be supported by NZD 2.8 billion in business tax changes, including an increase in the provisional tax threshold and the reintroduction of depreciation charges for commercial buildings. In addition, NZD 600 million will be spent on support for the aviation industry. In addition to the measures in the fiscal response package, the government put in place the Business Finance Guarantee Scheme (NZD 6.25 billion) for solvent SMEs, with the government bearing 80% of the risk of these loans and banks the remaining 20%. To support this scheme, the RBNZ introduced the Term Lending Facility (TLF) that provides funding for banks for terms of up to three years at low interest rates.
On 24 March, the government and RBNZ announced a support package for homeowners and businesses impacted by the economic effects of Covid-19. The package will include a six-month interest and principal payment holiday for mortgage holders and SMEs who have lost income because of the economic disruption caused by the Covid-19 crisis. The RBNZ will help banks to put this scheme in place with appropriate bank capital rules.

'@

$ws.Range("J43").Value = $j43

# --- NZL: Monetary policy / Macro-prudential regulation (K43) ---
$k43 = @'
On 15 March, the Reserve Bank of New Zealand (RBNZ) cut its policy rate (the overnight cash rate) by 75 basis points to 0.25%, and committed to maintain this rate for at least 12 months.
As regulatory relief, the RBNZ postponed the start date of the increased capital requirements for banks initially planned for 1 April by 12 months, to 1 July, 2021. In addition, the RBNZ has reduced banks’ core funding ratios (i.e., funding from deposits and long-term domestic borrowing) from 75% of assets to 50% to support increased lending.
On 20 March, the RBNZ announced measures to supply banks with more liquidity via both foreign exchange swaps and the reinstated Term Auction Facility, which offers banks term funding of up to one year against a range of collateral. It also put in place with the US Federal Reserve a USD 30 billion USD swap arrangement for at least six months.  
On 22 March, the RBNZ announced a Large Scale Asset Purchase programme (LSAP) of New Zealand government bonds (quantitative easing) amounting to NZD 30 billion (about 10% of GDP) over 12 months. The government will cover any loss the RBNZ incurs on the LSAP up until September 2021.
On 31 March, the RBNZ began weekly open market operations of up to NZD 500 million for corporate bonds and acceptable asset-backed securities, which is to be carried out over a year.

'@

$ws.Range("K43").Value = $k43

# Row 43 wraps less text now, so Excel recalculated a shorter row height
$ws.Rows.Item(43).RowHeight = 369.75

# --- TUR: Travel bans/restrictions (D54) ---
$d54 = @'
Air traffic is stopped with all countries. Land borders  with Iran, Azerbaijan and Georgia are closed. Public employees’ travels to foreign countries are subject to prior approval by their superiors. Turkish citizens are advised to postpone their travel plans abroad. Public transportation vehicles are required to accept 50 percent of their capacity to allow social distancing.                                               On 28 March, additional measures were announced: intercity travel was prohibited, subject only to individual permissions by state governors; Turkish Airlines suspended its domestic flights, and will only maintain some traffic between Istanbul, Ankara and a few other large cities.

'@

$ws.Range("D54").Value = $d54

# --- TUR: Health (G54) ---
$g54 = @'
The Ministry of Health (MoH) launched a country-wide public awereness campaign against Covid-19 and set up a call center for consultations to check prelimanry symptoms of corona before referring to hospitals, It increased the number of certified test centers and rapid diagnostic kits to raise the daily test capacity above 10.000 and started to recruit 32,000 additional health staff. Masks produced in Turkey will not be exported and ventilators are subject to export control. In-patient visits in hospitals are suspended during working hours and are limited to only one visitor after business hours. A periodic follow-up program including social services and home healthcare is introduced for seniors over 80 years of age and living alone. Protective masks and colognes are being distributed to all persons above 65 in Istanbul and Ankara. Comprehensive health surveillance of seniors residing in public nursing homes (more than 400 nursing homes, 153 out of which are public) is ensured and a "Coronavirus Information Guide for Nursing Homes" was distributed to all nursing homes. Preventive measures are  taken in women's shelters and in nursing homes for children under state protection. Daily disinfection activites are ongoing in public spaces and transportation facilities.
'@

$ws.Range("G54").Value = $g54

# --- Bump "Updated on" date for every existing data row (5-57) ---
$ws.Range("B5:B57").Value = 43924

# --- Append 17 new rows (58-74) that only carry the updated date ---
$ws.Range("B58:B74").Value = 43924
$ws.Range("B58:B74").NumberFormat = $ws.Range("B57").NumberFormat
$ws.Range("B58:B74").WrapText = $true
$ws.Range("B58:B74").VerticalAlignment = -4160

Write-Output "Done"
